# Slide 12 ("Assumptions") - reposition the picture + caption textbox, and
# bump the caption's font size to 24pt, matching the authored edit.
#
# NOTE: the source diff also wraps the slide transition in an
# mc:AlternateContent block offering a PowerPoint-2019+ "Morph" transition
# (p159:morph) with a fade fallback. Morph is a UI-only gallery transition -
# the PowerPoint object model has no SlideShowTransition member that can
# author it (touching .EntryEffect only lets us pick from the legacy
# PpEntryEffect gallery, e.g. fade/wipe/push/etc., and always stamps an
# unrelated p14:dur duration extension in the process), so that part of the
# edit cannot be reproduced through COM automation and is intentionally left
# as-is here.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)

# Locate the shapes by name so this doesn't depend on z-order assumptions.
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.Name -eq "Imagem 6") { $pic = $shape }
    if ($shape.Name -eq "CaixaDeTexto 2") { $tb = $shape }
}

# --- "Imagem 6" picture: move/resize (EMU -> points, 12700 EMU per point) ---
$pic.Left   = 229.65496072992124
$pic.Top    = 118.59228346456693
$pic.Width  = 453.54338582677167
$pic.Height = 340.1574860149606

# --- "CaixaDeTexto 2" caption textbox: move/resize ---
$tb.Left   = 75.5132283464567
$tb.Top    = 37.40992165984252
$tb.Width  = 404.4867716535433
$tb.Height = 36.3515759031496

# Bump every run in the caption to 24pt ("Importance" / " for " / "each" / " feature")
$tr = $tb.TextFrame.TextRange
$tr.Characters(1, 10).Font.Size = 24   # "Importance"
$tr.Characters(11, 5).Font.Size = 24   # " for "
$tr.Characters(16, 4).Font.Size = 24   # "each"
$tr.Characters(20, 8).Font.Size = 24   # " feature"
